$wb = $excel.ActiveWorkbook

# --- ETALONNAGE sheet: rows 27-31, columns B (num), C (text-numeric), D (num), E (num) ---
$wsEtalonnage = $wb.Worksheets.Item("ETALONNAGE")
$wsEtalonnage.Cells.Item(27, 2).Value = 5786417.22007224
$wsEtalonnage.Cells.Item(27, 3).Formula = "'153.828448933602"
$wsEtalonnage.Cells.Item(27, 4).Value = 4.56126308163687
$wsEtalonnage.Cells.Item(27, 5).Value = 8.00174227543069
$wsEtalonnage.Cells.Item(28, 2).Value = 5793686.17376681
$wsEtalonnage.Cells.Item(28, 3).Formula = "'154.089521210907"
$wsEtalonnage.Cells.Item(28, 4).Value = 0.125620974397589
$wsEtalonnage.Cells.Item(28, 5).Value = 0.169716511551643
$wsEtalonnage.Cells.Item(29, 2).Value = 6075527.15573691
$wsEtalonnage.Cells.Item(29, 3).Formula = "'164.212142717708"
$wsEtalonnage.Cells.Item(29, 4).Value = 4.86462285869473
$wsEtalonnage.Cells.Item(29, 5).Value = 6.56931206434603
$wsEtalonnage.Cells.Item(30, 2).Value = 6097442.06192298
$wsEtalonnage.Cells.Item(30, 3).Formula = "'164.99924010549"
$wsEtalonnage.Cells.Item(30, 4).Value = 0.360707896192736
$wsEtalonnage.Cells.Item(30, 5).Value = 0.479317409027136
$wsEtalonnage.Cells.Item(31, 2).Value = 6365418.3214698
$wsEtalonnage.Cells.Item(31, 3).Formula = "'174.623895218594"
$wsEtalonnage.Cells.Item(31, 4).Value = 4.39489636515387
$wsEtalonnage.Cells.Item(31, 5).Value = 5.83315117509056

# --- INDI sheet: rows 102-121, columns B, C (both numeric) ---
$wsIndi = $wb.Worksheets.Item("INDI")
$wsIndi.Cells.Item(102, 2).Value = 4027346.57628613
$wsIndi.Cells.Item(102, 3).Value = 146.267172768539
$wsIndi.Cells.Item(103, 2).Value = 4196038.32357095
$wsIndi.Cells.Item(103, 3).Value = 152.393803411659
$wsIndi.Cells.Item(104, 2).Value = 4183313.18378183
$wsIndi.Cells.Item(104, 3).Value = 151.931645466028
$wsIndi.Cells.Item(105, 2).Value = 4535462.35938847
$wsIndi.Cells.Item(105, 3).Value = 164.721174088184
$wsIndi.Cells.Item(106, 2).Value = 4332745.65644015
$wsIndi.Cells.Item(106, 3).Value = 157.358808209916
$wsIndi.Cells.Item(107, 2).Value = 4244432.86188409
$wsIndi.Cells.Item(107, 3).Value = 154.151420284808
$wsIndi.Cells.Item(108, 2).Value = 4327859.99230963
$wsIndi.Cells.Item(108, 3).Value = 157.181368234008
$wsIndi.Cells.Item(109, 2).Value = 4065875.57607889
$wsIndi.Cells.Item(109, 3).Value = 147.666488114894
$wsIndi.Cells.Item(110, 2).Value = 4329159.92558898
$wsIndi.Cells.Item(110, 3).Value = 157.22857985634
$wsIndi.Cells.Item(111, 2).Value = 4484142.78805569
$wsIndi.Cells.Item(111, 3).Value = 162.857324413379
$wsIndi.Cells.Item(112, 2).Value = 4476546.50858915
$wsIndi.Cells.Item(112, 3).Value = 162.58143896372
$wsIndi.Cells.Item(113, 2).Value = 4795937.17100717
$wsIndi.Cells.Item(113, 3).Value = 174.181227637392
$wsIndi.Cells.Item(114, 2).Value = 4618077.44026453
$wsIndi.Cells.Item(114, 3).Value = 167.721629618616
$wsIndi.Cells.Item(115, 2).Value = 4542742.28460536
$wsIndi.Cells.Item(115, 3).Value = 164.985570027116
$wsIndi.Cells.Item(116, 2).Value = 4621314.13335757
$wsIndi.Cells.Item(116, 3).Value = 167.839181445575
$wsIndi.Cells.Item(117, 2).Value = 4390340.85775581
$wsIndi.Cells.Item(117, 3).Value = 159.450579330653
$wsIndi.Cells.Item(118, 2).Value = 4630094.1313516
$wsIndi.Cells.Item(118, 3).Value = 168.1580577725
$wsIndi.Cells.Item(119, 2).Value = 4772791.59543817
$wsIndi.Cells.Item(119, 3).Value = 173.340615130757
$wsIndi.Cells.Item(120, 2).Value = 4769791.61196034
$wsIndi.Cells.Item(120, 3).Value = 173.231660241144
$wsIndi.Cells.Item(121, 2).Value = 5059825.29967156
$wsIndi.Cells.Item(121, 3).Value = 183.765247729975

# --- PREVISION sheet: rows 98-117, columns B,C,D,E,F,I,J (numeric), K (text) ---
$wsPrevision = $wb.Worksheets.Item("PREVISION")
$wsPrevision.Cells.Item(98, 2).Value = 1393972.74367111
$wsPrevision.Cells.Item(98, 3).Value = 690624.228515737
$wsPrevision.Cells.Item(98, 4).Value = 692680.858133599
$wsPrevision.Cells.Item(98, 5).Value = 701291.885537515
$wsPrevision.Cells.Item(98, 6).Value = 36.5667931921348
$wsPrevision.Cells.Item(98, 9).Value = 1.84389360174437
$wsPrevision.Cells.Item(98, 10).Value = 2.14717744168413
$wsPrevision.Cells.Item(98, 11).Value = "Acceptable"
$wsPrevision.Cells.Item(99, 2).Value = 1436618.20847153
$wsPrevision.Cells.Item(99, 3).Value = 733497.920569989
$wsPrevision.Cells.Item(99, 4).Value = 713871.873013605
$wsPrevision.Cells.Item(99, 5).Value = 722746.335457923
$wsPrevision.Cells.Item(99, 6).Value = 38.0984508529146
$wsPrevision.Cells.Item(99, 9).Value = 6.99131837721334
$wsPrevision.Cells.Item(99, 10).Value = 4.12857447064627
$wsPrevision.Cells.Item(100, 2).Value = 1433401.27871834
$wsPrevision.Cells.Item(100, 3).Value = 727391.262182339
$wsPrevision.Cells.Item(100, 4).Value = 712273.344152758
$wsPrevision.Cells.Item(100, 5).Value = 721127.934565585
$wsPrevision.Cells.Item(100, 6).Value = 37.9829113665069
$wsPrevision.Cells.Item(100, 9).Value = 0.640617661811027
$wsPrevision.Cells.Item(100, 10).Value = -1.45107175949549
$wsPrevision.Cells.Item(101, 2).Value = 1522424.98921126
$wsPrevision.Cells.Item(101, 3).Value = 823511.069561726
$wsPrevision.Cells.Item(101, 4).Value = 756510.235052125
$wsPrevision.Cells.Item(101, 5).Value = 765914.754159133
$wsPrevision.Cells.Item(101, 6).Value = 41.1802935220461
$wsPrevision.Cells.Item(101, 9).Value = 24.1243253725232
$wsPrevision.Cells.Item(101, 10).Value = 14.0255741956552
$wsPrevision.Cells.Item(102, 2).Value = 1471177.97707709
$wsPrevision.Cells.Item(102, 3).Value = 742995.187936689
$wsPrevision.Cells.Item(102, 4).Value = 731045.013796514
$wsPrevision.Cells.Item(102, 5).Value = 740132.96328058
$wsPrevision.Cells.Item(102, 6).Value = 39.3397020524789
$wsPrevision.Cells.Item(102, 9).Value = 7.58313381699702
$wsPrevision.Cells.Item(102, 10).Value = 5.53850380192193
$wsPrevision.Cells.Item(103, 2).Value = 1448852.40257636
$wsPrevision.Cells.Item(103, 3).Value = 741957.636731357
$wsPrevision.Cells.Item(103, 4).Value = 719951.182748736
$wsPrevision.Cells.Item(103, 5).Value = 728901.219827622
$wsPrevision.Cells.Item(103, 6).Value = 38.537855071202
$wsPrevision.Cells.Item(103, 9).Value = 1.15333880630426
$wsPrevision.Cells.Item(103, 10).Value = 0.851596759158868
$wsPrevision.Cells.Item(104, 2).Value = 1469942.87565063
$wsPrevision.Cells.Item(104, 3).Value = 752524.949496759
$wsPrevision.Cells.Item(104, 4).Value = 730431.277896833
$wsPrevision.Cells.Item(104, 5).Value = 739511.597753793
$wsPrevision.Cells.Item(104, 6).Value = 39.295342058502
$wsPrevision.Cells.Item(104, 9).Value = 3.45531883886179
$wsPrevision.Cells.Item(104, 10).Value = 2.54929289340082
$wsPrevision.Cells.Item(105, 2).Value = 1403712.91846273
$wsPrevision.Cells.Item(105, 3).Value = 738247.455064996
$wsPrevision.Cells.Item(105, 4).Value = 697520.861400277
$wsPrevision.Cells.Item(105, 5).Value = 706192.057062453
$wsPrevision.Cells.Item(105, 6).Value = 36.9166220287235
$wsPrevision.Cells.Item(105, 9).Value = -10.3536695070908
$wsPrevision.Cells.Item(105, 10).Value = -7.79756451646464
$wsPrevision.Cells.Item(106, 2).Value = 1470271.50025615
$wsPrevision.Cells.Item(106, 3).Value = 742380.293599723
$wsPrevision.Cells.Item(106, 4).Value = 730594.575188544
$wsPrevision.Cells.Item(106, 5).Value = 739676.92506761
$wsPrevision.Cells.Item(106, 6).Value = 39.3071449640849
$wsPrevision.Cells.Item(106, 9).Value = -0.0827588585969208
$wsPrevision.Cells.Item(106, 10).Value = -0.0616157144182727
$wsPrevision.Cells.Item(107, 2).Value = 1509451.34344571
$wsPrevision.Cells.Item(107, 3).Value = 783860.669742079
$wsPrevision.Cells.Item(107, 4).Value = 750063.483404503
$wsPrevision.Cells.Item(107, 5).Value = 759387.860041207
$wsPrevision.Cells.Item(107, 6).Value = 40.7143311033447
$wsPrevision.Cells.Item(107, 9).Value = 5.64763147331744
$wsPrevision.Cells.Item(107, 10).Value = 4.18254756396126
$wsPrevision.Cells.Item(108, 2).Value = 1507530.99539183
$wsPrevision.Cells.Item(108, 3).Value = 778378.445994548
$wsPrevision.Cells.Item(108, 4).Value = 749109.240687838
$wsPrevision.Cells.Item(108, 5).Value = 758421.75470399
$wsPrevision.Cells.Item(108, 6).Value = 40.6453597409301
$wsPrevision.Cells.Item(108, 9).Value = 3.43556668985896
$wsPrevision.Cells.Item(108, 10).Value = 2.55711431810337
$wsPrevision.Cells.Item(109, 2).Value = 1588273.31664321
$wsPrevision.Cells.Item(109, 3).Value = 870805.892826211
$wsPrevision.Cells.Item(109, 4).Value = 789231.015396873
$wsPrevision.Cells.Item(109, 5).Value = 799042.30124634
$wsPrevision.Cells.Item(109, 6).Value = 43.545306909348
$wsPrevision.Cells.Item(109, 9).Value = 17.9558272570739
$wsPrevision.Cells.Item(109, 10).Value = 13.1480159335289
$wsPrevision.Cells.Item(110, 2).Value = 1543310.17523956
$wsPrevision.Cells.Item(110, 3).Value = 791924.93345078
$wsPrevision.Cells.Item(110, 4).Value = 766888.320739988
$wsPrevision.Cells.Item(110, 5).Value = 776421.854499573
$wsPrevision.Cells.Item(110, 6).Value = 41.9304074046541
$wsPrevision.Cells.Item(110, 9).Value = 6.67375471550047
$wsPrevision.Cells.Item(110, 10).Value = 4.96769984119823
$wsPrevision.Cells.Item(111, 2).Value = 1524265.36255247
$wsPrevision.Cells.Item(111, 3).Value = 794104.286589947
$wsPrevision.Cells.Item(111, 4).Value = 757424.737427489
$wsPrevision.Cells.Item(111, 5).Value = 766840.625124986
$wsPrevision.Cells.Item(111, 6).Value = 41.246392506779
$wsPrevision.Cells.Item(111, 9).Value = 1.30681602525591
$wsPrevision.Cells.Item(111, 10).Value = 0.981417464768874
$wsPrevision.Cells.Item(112, 2).Value = 1544128.41491987
$wsPrevision.Cells.Item(112, 3).Value = 803550.528666168
$wsPrevision.Cells.Item(112, 4).Value = 767294.913312539
$wsPrevision.Cells.Item(112, 5).Value = 776833.501607331
$wsPrevision.Cells.Item(112, 6).Value = 41.9597953613938
$wsPrevision.Cells.Item(112, 9).Value = 3.23391311786099
$wsPrevision.Cells.Item(112, 10).Value = 2.42763960674179
$wsPrevision.Cells.Item(113, 2).Value = 1485738.10921108
$wsPrevision.Cells.Item(113, 3).Value = 797161.129124334
$wsPrevision.Cells.Item(113, 4).Value = 738280.108504711
$wsPrevision.Cells.Item(113, 5).Value = 747458.000706366
$wsPrevision.Cells.Item(113, 6).Value = 39.8626448326633
$wsPrevision.Cells.Item(113, 9).Value = -8.45708145851683
$wsPrevision.Cells.Item(113, 10).Value = -6.45576591621162
$wsPrevision.Cells.Item(114, 2).Value = 1546348.00835837
$wsPrevision.Cells.Item(114, 3).Value = 793985.599910476
$wsPrevision.Cells.Item(114, 4).Value = 768397.854453008
$wsPrevision.Cells.Item(114, 5).Value = 777950.153905361
$wsPrevision.Cells.Item(114, 6).Value = 42.0395144431251
$wsPrevision.Cells.Item(114, 9).Value = 0.260209821998747
$wsPrevision.Cells.Item(114, 10).Value = 0.196838792845799
$wsPrevision.Cells.Item(115, 2).Value = 1582422.08892105
$wsPrevision.Cells.Item(115, 3).Value = 834318.573999225
$wsPrevision.Cells.Item(115, 4).Value = 786323.474013352
$wsPrevision.Cells.Item(115, 5).Value = 796098.614907696
$wsPrevision.Cells.Item(115, 6).Value = 43.3351537826893
$wsPrevision.Cells.Item(115, 9).Value = 5.06410657748317
$wsPrevision.Cells.Item(115, 10).Value = 3.81539381510234
$wsPrevision.Cells.Item(116, 2).Value = 1581663.6896996
$wsPrevision.Cells.Item(116, 3).Value = 829367.677854335
$wsPrevision.Cells.Item(116, 4).Value = 785946.616843149
$wsPrevision.Cells.Item(116, 5).Value = 795717.072856452
$wsPrevision.Cells.Item(116, 6).Value = 43.3079150602859
$wsPrevision.Cells.Item(116, 9).Value = 3.21288435103411
$wsPrevision.Cells.Item(116, 10).Value = 2.4308389391098
$wsPrevision.Cells.Item(117, 2).Value = 1654984.53449078
$wsPrevision.Cells.Item(117, 3).Value = 918720.477461934
$wsPrevision.Cells.Item(117, 4).Value = 822380.575770697
$wsPrevision.Cells.Item(117, 5).Value = 832603.958720083
$wsPrevision.Cells.Item(117, 6).Value = 45.9413119324937
$wsPrevision.Cells.Item(117, 9).Value = 15.2490310799688
$wsPrevision.Cells.Item(117, 10).Value = 11.3914036552224

# --- VATRIM sheet: rows 98-117, column B (numeric) ---
$wsVatrim = $wb.Worksheets.Item("VATRIM")
$wsVatrim.Cells.Item(98, 2).Value = 690624.228515737
$wsVatrim.Cells.Item(99, 2).Value = 733497.920569989
$wsVatrim.Cells.Item(100, 2).Value = 727391.262182339
$wsVatrim.Cells.Item(101, 2).Value = 823511.069561726
$wsVatrim.Cells.Item(102, 2).Value = 742995.187936689
$wsVatrim.Cells.Item(103, 2).Value = 741957.636731357
$wsVatrim.Cells.Item(104, 2).Value = 752524.949496759
$wsVatrim.Cells.Item(105, 2).Value = 738247.455064996
$wsVatrim.Cells.Item(106, 2).Value = 742380.293599723
$wsVatrim.Cells.Item(107, 2).Value = 783860.669742079
$wsVatrim.Cells.Item(108, 2).Value = 778378.445994548
$wsVatrim.Cells.Item(109, 2).Value = 870805.892826211
$wsVatrim.Cells.Item(110, 2).Value = 791924.93345078
$wsVatrim.Cells.Item(111, 2).Value = 794104.286589947
$wsVatrim.Cells.Item(112, 2).Value = 803550.528666168
$wsVatrim.Cells.Item(113, 2).Value = 797161.129124334
$wsVatrim.Cells.Item(114, 2).Value = 793985.599910476
$wsVatrim.Cells.Item(115, 2).Value = 834318.573999225
$wsVatrim.Cells.Item(116, 2).Value = 829367.677854335
$wsVatrim.Cells.Item(117, 2).Value = 918720.477461934
